# repull data, push all data, mean calculation
# Update the dSF column (F) values for the re-pulled rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = 1
$ws.Range("F5").Value = 2
$ws.Range("F7").Value = -1
$ws.Range("F9").Value = 2
